$d = $word.ActiveDocument

# Common run-properties fragments reused while building OOXML snippets.
$rprPlain = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'
$rprBody  = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# Step 1: locate the "Tuan 6 ... Trang login va logout" heading paragraph and
# collapse its 4 runs ("Tuan " / "6" / " " / "Trang login va logout") down to
# the 2 runs the edit leaves behind ("Tuan " / "6 Trang login va logout").
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute("Trang login và logout", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "could not find the week-6 heading paragraph"
}
$headingPara = $d.Content.Find.Parent.Paragraphs.Item(1)

# Re-resolve the paragraph via the Paragraphs collection using the found range.
$p6 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Trang login và logout*") {
        $p6 = $cand
        break
    }
}
if ($null -eq $p6) {
    throw "could not resolve week-6 heading paragraph"
}

$p6Start = $p6.Range.Start
$p6End = $p6.Range.End - 1
$r6 = $d.Range($p6Start, $p6End)

$week6Xml = $pkgHeader + '<w:body>' `
  + '<w:p>' `
  +   '<w:r>' + $rprPlain + '<w:t xml:space="preserve">Tuần </w:t></w:r>' `
  +   '<w:r>' + $rprPlain + '<w:t>6 Trang login và logout</w:t></w:r>' `
  + '</w:p>' `
  + '</w:body>' + $pkgFooter
$r6.InsertXML($week6Xml)

# ---------------------------------------------------------------------------
# Step 2: right after the "Biet cach dung Spring Security..." paragraph
# (which stays put) insert the new "Tuan 7" block: heading + 3 body
# paragraphs + a trailing empty paragraph.
# ---------------------------------------------------------------------------
$pSpring = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Biết cách dùng Spring Security*") {
        $pSpring = $cand
        break
    }
}
if ($null -eq $pSpring) {
    throw "could not resolve Spring Security paragraph"
}

$insertPoint = $d.Range($pSpring.Range.End, $pSpring.Range.End)

$week7Body = '<w:body>' `
  + '<w:p>' `
  +   '<w:pPr><w:pStyle w:val="Heading2"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' + $rprPlain + '</w:pPr>' `
  +   '<w:r>' + $rprPlain + '<w:t xml:space="preserve">Tuần </w:t></w:r>' `
  +   '<w:r>' + $rprPlain + '<w:t>7</w:t></w:r>' `
  +   '<w:r>' + $rprPlain + '<w:t xml:space="preserve"> </w:t></w:r>' `
  +   '<w:r>' + $rprPlain + '<w:t>Kết nối cơ sở dữ liệu</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p>' `
  +   '<w:pPr>' + $rprBody + '</w:pPr>' `
  +   '<w:r>' + $rprBody + '<w:t>Tạo cở sở dữ liệu để kết nối với project</w:t></w:r>' `
  +   '<w:r>' + $rprBody + '<w:t>.</w:t></w:r>' `
  +   '<w:r>' + $rprBody + '<w:t xml:space="preserve"> Dùng JPA để kết nối MYSQL.</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p>' `
  +   '<w:pPr>' + $rprBody + '</w:pPr>' `
  +   '<w:r>' + $rprBody + '<w:t>Thiết kế cở sở dữ liệu để load lên category và product dưới sự quản lý của admin.</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p>' `
  +   '<w:pPr>' + $rprBody + '</w:pPr>' `
  + '</w:p>' `
  + '</w:body>'

$week7Xml = $pkgHeader + $week7Body + $pkgFooter
$insertPoint.InsertXML($week7Xml)

# ---------------------------------------------------------------------------
# Step 3: append the new progress paragraph right after the "... 21/11/2025
# ..." summary paragraph.
# ---------------------------------------------------------------------------
$pProgress = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*21/11/2025*") {
        $pProgress = $cand
        break
    }
}
if ($null -eq $pProgress) {
    throw "could not resolve 21/11/2025 progress paragraph"
}

$progressInsertPoint = $d.Range($pProgress.Range.End, $pProgress.Range.End)

$progressXml = $pkgHeader + '<w:body>' `
  + '<w:p>' `
  +   '<w:pPr>' + $rprBody + '</w:pPr>' `
  +   '<w:r>' + $rprBody + '<w:t>Tỷ lệ hoàn thành đến ngày 27/11/2025: Hoàn thành kết nối cơ sở dữ liệu lên project, bước đầu thành công với 2 trang category và product.</w:t></w:r>' `
  + '</w:p>' `
  + '</w:body>' + $pkgFooter
$progressInsertPoint.InsertXML($progressXml)

Write-Output "done"
